$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Vsds")
$ws.Rows.Item(22).Delete()
